$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = 900
$ws.Range("B17").Value = '71.473.820/0012-12'
$ws.Range("C17").Value = 'MILLENNIUM CCVM S/A, SUCESSORA DA GAMEX SECURITIES CCVM LTDA.'
$ws.Range("D17").Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=IA&Ano=2002&NumProc=15'

$ws.Range("A18").Value = 1456
$ws.Range("B18").Value = '47.894.290/0001-28'
$ws.Range("C18").Value = 'PLANIN AUDITORES INDEPENDENTES  S/C'
$ws.Range("D18").Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2001&NumProc=8739'

$ws.Range("A19").Value = 397
$ws.Range("B19").Value = '09.143.363/0001-50'
$ws.Range("C19").Value = 'EASE ESCRITÓRIO DE AUDITORIA INDEPENDENTE S/C'
$ws.Range("D19").Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2013&NumProc=13481'

$ws.Range("A20").Value = 411
$ws.Range("B20").Value = '62.030.762/0001-98'
$ws.Range("C20").Value = 'AKW AUDITORES INDEPENDENTES S/S'
$ws.Range("D20").Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2013&NumProc=4362'

$ws.Range("A21").Value = 1165
$ws.Range("B21").Value = '67.634.717/0001-66'
$ws.Range("C21").Value = 'BWEL AUDITORES INDEPENDENTES SOCIEDADE SIMPLES'
$ws.Range("D21").Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2017&NumProc=2253'

$ws.Range("A22").Value = 427
$ws.Range("B22").Value = '11.245.719/0003-70'
$ws.Range("C22").Value = 'DIRECTA AUDITORES'
$ws.Range("D22").Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2013&NumProc=5682'

$ws.Range("A23").Value = 1688
$ws.Range("B23").Value = '58.214.958/0001-65'
$ws.Range("C23").Value = 'SOC TEC AUDITORIA SOMATEC SC'
$ws.Range("D23").Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=IA&Ano=2000&NumProc=6'

$ws.Range("A24").Value = 844
$ws.Range("B24").Value = '31.622.483/0001-90'
$ws.Range("C24").Value = 'GUILHERME FONTES FILMES LTDA.'
$ws.Range("D24").Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=IA&Ano=2000&NumProc=12'

$ws.Range("A25").Value = 1035
$ws.Range("B25").Value = '00.469.585/0001-93'
$ws.Range("C25").Value = 'FACEB - FUNDAÇÃO DE ASSISTÊNCIA DOS EMPREGADOS DA CEB'
$ws.Range("D25").Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=IA&Ano=1999&NumProc=28'

$ws.Range("A26").Value = 1035
$ws.Range("B26").Value = '17.393.471/0001-13'
$ws.Range("C26").Value = 'PRATA DTVM LTDA. (atual Prata Consultoria e Assessoria Empresarial Ltda.)'
$ws.Range("D26").Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=IA&Ano=1999&NumProc=28'

$ws.Range("A27").Value = 52
$ws.Range("B27").Value = '04.612.682/0001-44'
$ws.Range("C27").Value = 'INTERTRADING AGRONEGÓCIOS LTDA.'
$ws.Range("D27").Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2007&NumProc=4414'

$ws.Range("A28").Value = 21
$ws.Range("B28").Value = '27.901.719/0001-50'
$ws.Range("C28").Value = 'INSTITUTO AERUS DE SEGURIDADE SOCIAL'
$ws.Range("D28").Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2007&NumProc=1176'

$ws.Range("A29").Value = 1563
$ws.Range("B29").Value = '43.214.485/0001-29'
$ws.Range("C29").Value = 'SOCIVAL AUDITORIA INDEP SC'
$ws.Range("D29").Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2004&NumProc=7001'

$ws.Range("A30").Value = 555
$ws.Range("B30").Value = '64.920.416/0001-00'
$ws.Range("C30").Value = 'NORMAS AUDITORES INDEPENDENTES'
$ws.Range("D30").Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2015&NumProc=11941'

$ws.Range("A31").Value = 1565
$ws.Range("B31").Value = '43.729.789/0001-29'
$ws.Range("C31").Value = 'PERMALI DO BRASIL INDÚSTRIA E COMÉRCIO LTDA'
$ws.Range("D31").Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2005&NumProc=33'

$ws.Range("A32").Value = 1556
$ws.Range("B32").Value = '04.565.230/0002-30'
$ws.Range("C32").Value = 'I.B. Sabbá S/A'
$ws.Range("D32").Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2004&NumProc=4627'

$ws.Range("A33").Value = 217
$ws.Range("B33").Value = '05.723.617/0001-59'
$ws.Range("C33").Value = 'MAPFRE DTVM S.A.'
$ws.Range("D33").Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2010&NumProc=17292'

$ws.Range("A34").Value = 1710
$ws.Range("B34").Value = '14.629.882/0001-63'
$ws.Range("C34").Value = 'CAPITAL ASSESSORIA FINANCEIRA LTDA. (ATUAL CAPITAL ASSESSORIA E EMPREENDIMENTOS LTDA.)'
$ws.Range("D34").Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=IA&Ano=2002&NumProc=6'

$ws.Range("A35").Value = 642
$ws.Range("B35").Value = '05.706.592/0001-85'
$ws.Range("C35").Value = 'BANCO BOZANO, SIMONSEN S/A'
$ws.Range("D35").Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=IA&Ano=1999&NumProc=10'

$ws.Range("A36").Value = 2
$ws.Range("B36").Value = '00.659.559/0002-09'
$ws.Range("C36").Value = 'MASTER CORRETORA DE MERCADORIAS LTDA.'
$ws.Range("D36").Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=IA&Ano=2006&NumProc=1'

$ws.Range("A37").Value = 4
$ws.Range("B37").Value = '00.016.087/6747-72'
$ws.Range("C37").Value = 'MARCOS LEVY'
$ws.Range("D37").Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2006&NumProc=8625'
